# New weekly Orégano price point (Mercado Mayorista Lo Valledor de Santiago).
# Insert one new data row at row 263 (pushing existing rows 263:280 down to
# 264:281) and fill it with the latest reading, carrying forward the same
# volume/price figures as the most recent existing record but dated later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("263:263").Insert()

$ws.Range("A263").Value = 6
$ws.Range("B263").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C263").Value = "Metropolitana"
$ws.Range("D263").Value = 44931
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 100112029
$ws.Range("G263").Value = "Orégano"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 46
$ws.Range("K263").Value = 16000
$ws.Range("L263").Value = 17000
$ws.Range("M263").Value = 16457
$ws.Range("N263").Value = "$/docena de atados"
$ws.Range("O263").Value = "Región Metropolitana"
$ws.Range("P263").Value = 5486
$ws.Range("Q263").Value = 3
$ws.Range("R263").Value = "Hortaliza"
